$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New RF values (old value / 1.5, rounded to 2 decimals) for columns E, I, M across rows 2-10
$values = @{
    "E2" = 0.41; "I2" = 0.6;  "M2" = 0.82
    "E3" = 0.36; "I3" = 0.55; "M3" = 0.72
    "E4" = 0.36; "I4" = 0.55; "M4" = 0.73
    "E5" = 0.37; "I5" = 0.55; "M5" = 0.73
    "E6" = 0.37; "I6" = 0.55; "M6" = 0.73
    "E7" = 0.37; "I7" = 0.55; "M7" = 0.73
    "E8" = 0.36; "I8" = 0.55; "M8" = 0.73
    "E9" = 0.36; "I9" = 0.55; "M9" = 0.72
    "E10" = 0.41; "I10" = 0.6;  "M10" = 0.82
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
